{"js": "// Remove the empty paragraph, the \"Ver no Jupiter...\" paragraph, and the\n// \"\u00a9 2020 ...\" copyright paragraph that immediately follow the\n// \"LOQ4086: Opera\u00e7\u00f5es Unit\u00e1rias II (Requisito fraco)\" paragraph.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\nconst marker = \"LOQ4086: Opera\u00e7\u00f5es Unit\u00e1rias II (Requisito fraco)\";\nconst ver = \"Ver no Jupiter Salvar em pdf Salvar em docx\";\nconst copyright = \"Contact: luizeleno@usp.br\";\n\nlet markerIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text === marker) {\n    markerIndex = i;\n    break;\n  }\n}\n\nif (markerIndex === -1) {\n  throw new Error(\"Could not locate the LOQ4086 requirement paragraph.\");\n}\n\n// The three paragraphs to remove are the ones right after the marker:\n// an empty paragraph, the \"Ver no Jupiter...\" paragraph, and the\n// copyright paragraph. Confirm via content before deleting, so the\n// script is resilient to being run only when the pattern truly matches.\nconst toDelete = [];\nfor (let i = markerIndex + 1; i < items.length && toDelete.length < 3; i++) {\n  const t = items[i].text;\n  if (t === \"\" || t === ver || t.indexOf(copyright) !== -1) {\n    toDelete.push(items[i]);\n  } else {\n    break;\n  }\n}\n\nfor (const p of toDelete) {\n  p.delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the empty paragraph, the \"Ver no Jupiter...\" paragraph, and the\n# \"\u00a9 2020 ...\" copyright paragraph that immediately follow the\n# \"LOQ4086: Opera\u00e7\u00f5es Unit\u00e1rias II (Requisito fraco)\" paragraph.\n\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$found = $rng.Find.Execute(\"LOQ4086: Opera\u00e7\u00f5es Unit\u00e1rias II (Requisito fraco)\")\nif (-not $found) {\n    throw \"Could not locate the LOQ4086 requirement paragraph.\"\n}\n\n# Resolve which paragraph (by 1-based index) the found range falls inside.\n$count = $d.Paragraphs.Count\n$markerIndex = -1\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Start -le $rng.Start -and $p.Range.End -ge $rng.End) {\n        $markerIndex = $i\n        break\n    }\n}\nif ($markerIndex -eq -1) {\n    throw \"Could not resolve the paragraph index of the LOQ4086 requirement.\"\n}\n\n# The three paragraphs to delete are the ones immediately following the\n# marker paragraph: an empty paragraph, \"Ver no Jupiter...\", and the\n# copyright notice.\n$firstToDelete = $d.Paragraphs.Item($markerIndex + 1)\n$lastToDelete = $d.Paragraphs.Item($markerIndex + 3)\n\n$delRange = $d.Range($firstToDelete.Range.Start, $lastToDelete.Range.End)\n$delRange.Delete()\n"}
